$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace Chert_type values: "A" -> "Coarser", "B" -> "Finer" in column B (rows 2-13)
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $current = $cell.Value2
    if ($current -eq "A") {
        $cell.Value = "Coarser"
    } elseif ($current -eq "B") {
        $cell.Value = "Finer"
    }
}
